$wb = $excel.ActiveWorkbook

# Rename the "SwateTemplateMetadata" sheet to "isa_template"
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Name = "isa_template"
